# Generate Report for Handoff
# Replace old GUID-based file identifiers / hashes and timestamps with the
# newly generated ones, on all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "303dea24-92f5-4b42-865d-0c761ff8ec15"
$newGuid = "2d391c22-86b7-4001-8730-e8fa8ae5e9ac"

$oldHash = "e1ae066e8a2cd6ac3f0a281f876cd4f469b4abf9"
$newHash = "95fac6b7bc3bef435c0657a266a9a396aaa81722"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-31 04:58:59"

# ---- zh-cn sheet ----
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-31 04:58:53"

# ---- de-de sheet ----
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-31 04:58:59"

# ---- Update the displayed text of the hyperlinks so it matches the new
#      file names. The underlying hyperlink target itself is untouched by
#      this change (only the friendly display text is regenerated), so we
#      keep pointing at the same "blob" URL that was already present. The
#      simulated Hyperlinks collection does not support in-place mutation
#      of a Hyperlink object's properties, so the collection is rebuilt
#      (Delete + Add) using that same target address.

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6181744214685d904365602af9f15ee6ea2347e3/e2e/$oldGuid.md"

# Original hyperlink font: single underline, RGB FF6495ED (as an OLE/VBA
# BGR-packed color value).
$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

function Update-HyperlinkDisplay($ws, $cellAddr, $display, $targetUrl) {
    $range = $ws.Range($cellAddr)

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($range, $targetUrl, "", "", $display)

    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

Update-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md" $hyperlinkTarget
Update-HyperlinkDisplay $wsZhCn "A2" "$newGuid.md" $hyperlinkTarget
Update-HyperlinkDisplay $wsDeDe "A2" "$newGuid.md" $hyperlinkTarget

Write-Host "Localization status report regenerated for handoff."
